$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (H1) onto the two new
# header cells so they pick up the same style (bold, bordered, centered)
# used by the rest of row 1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I (I0) and J (IF) columns, rows 2-29.
$data = @{
    2  = @(9, 9)
    3  = @(9, 9)
    4  = @(10, 10)
    5  = @(9, 9)
    6  = @(7, 7)
    7  = @(8, 9)
    8  = @(5, 6)
    9  = @(5, 5)
    10 = @(8, 8)
    11 = @(3, 4)
    12 = @(7, 7)
    13 = @(5, 5)
    14 = @(8, 8)
    15 = @(7, 7)
    16 = @(6, 6)
    17 = @(5, 6)
    18 = @(9, 9)
    19 = @(5, 5)
    20 = @(6, 7)
    21 = @(10, 10)
    22 = @(9, 9)
    23 = @(8, 9)
    24 = @(9, 9)
    25 = @(8, 8)
    26 = @(6, 6)
    27 = @(3, 3)
    28 = @(3, 4)
    29 = @(7, 7)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
